# Added new code for pwm-ing an rgb led -- workbook side: add a "TDSez"
# column (new column B) ahead of the existing AtTiny / Atlas Sci columns,
# and fill in the calibration data rows below the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing header cells right by one column (C1 <- old B1 "AtTiny",
# D1 <- old C1 "Atlas Sci") before writing the new B1 header, so we don't
# clobber values we still need to move.
$ws.Range("D1").Value = $ws.Range("C1").Value2
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("B1").Value = "TDSez"

# Calibration data rows (A:D), row 2 through row 11.
$rows = @(
    @(1, 60, 145, 158),
    @(2, 110, 259, 244),
    @(3, 150, 547, 403),
    @(4, 183, 797, 661),
    @(5, 325, 1250, 1455),
    @(6, 429, 1531, 1778),
    @(7, 521, 1695, 2019),
    @(8, 739, 2039, 2471),
    @(9, 1070, 2970, 2906),
    @(10, 863, 2722, 2624)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Trailing summary rows with text labels in column A.
$ws.Cells.Item(12, 1).Value = "cal1 (220)"
$ws.Cells.Item(12, 2).Value = 114
$ws.Cells.Item(12, 3).Value = 220
$ws.Cells.Item(12, 4).Value = 212

$ws.Cells.Item(13, 1).Value = "cal2 (3000)"
$ws.Cells.Item(13, 2).Value = 1110
$ws.Cells.Item(13, 3).Value = 3000
$ws.Cells.Item(13, 4).Value = 3100

# Match the saved selection/active cell.
$ws.Range("G9").Select() | Out-Null
